$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '70.767.29'
$ws.Range('E2').Value = '  -2.04%  '
$ws.Range('D3').Value = '3.938.57'
$ws.Range('E3').Value = '  -2.34%  '
$ws.Range('E4').Value = '  -0.03%  '
$ws.Range('D5').Value = '536.22'
$ws.Range('E5').Value = '  +3.04%  '
$ws.Range('D6').Value = '147.96'
$ws.Range('E6').Value = '  +0.64%  '
$ws.Range('D7').Value = '3.935.59'
$ws.Range('E7').Value = '  -2.25%  '
$ws.Range('D8').Value = '0.685'
$ws.Range('E8').Value = '  -6.56%  '
$ws.Range('E9').Value = '  +0.03%  '
$ws.Range('D10').Value = '0.736'
$ws.Range('E10').Value = '  -5.32%  '
$ws.Range('D11').Value = '0.165'
$ws.Range('E11').Value = '  -6.13%  '
$ws.Range('D12').Value = '55.56'
$ws.Range('E12').Value = '  +15.10%  '
$ws.Range('D13').Value = '0.0000314'
$ws.Range('E13').Value = '  -4.07%  '
$ws.Range('D14').Value = '10.59'
$ws.Range('E14').Value = '  -4.68%  '
$ws.Range('D15').Value = '4.588.04'
$ws.Range('E15').Value = '  -1.93%  '
$ws.Range('D16').Value = '3.954.96'
$ws.Range('E16').Value = '  -2.43%  '
$ws.Range('D17').Value = '20.49'
$ws.Range('E17').Value = '  -3.44%  '
$ws.Range('D18').Value = '13.77'
$ws.Range('E18').Value = '  -2.79%  '
$ws.Range('E19').Value = '  -1.48%  '
$ws.Range('D20').Value = '1.16'
$ws.Range('E20').Value = '  -4.38%  '
$ws.Range('D21').Value = '70.765.38'
$ws.Range('E21').Value = '  -1.95%  '
$ws.Range('D22').Value = '424.15'
$ws.Range('E22').Value = '  -4.59%  '
$ws.Range('B23').Value = 'Litecoin'
$ws.Range('C23').Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range('D23').Value = '96.79'
$ws.Range('E23').Value = '  -8.04%  '
$ws.Range('B24').Value = 'ImmutableX'
$ws.Range('C24').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D24').Value = '3.56'
$ws.Range('E24').Value = '  -0.23%  '
$ws.Range('D25').Value = '4.21'
$ws.Range('E25').Value = '  +4.93%  '
$ws.Range('D26').Value = '14.37'
$ws.Range('E26').Value = '  -4.07%  '
$ws.Range('D27').Value = '11.30'
$ws.Range('E27').Value = '  -1.73%  '
$ws.Range('B28').Value = 'Filecoin'
$ws.Range('C28').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D28').Value = '10.60'
$ws.Range('E28').Value = '  -3.87%  '
$ws.Range('B29').Value = 'Toncoin'
$ws.Range('C29').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('D29').Value = '3.75'
$ws.Range('E29').Value = '  +15.45%  '
$ws.Range('E30').Value = '  +1.53%  '
$ws.Range('D31').Value = '36.31'
$ws.Range('E31').Value = '  -3.94%  '
$ws.Range('D32').Value = '7.73'
$ws.Range('E32').Value = '  +14.57%  '
$ws.Range('D33').Value = '49.93'
$ws.Range('E33').Value = '  +17.85%  '
$ws.Range('B34').Value = 'Hedera'
$ws.Range('C34').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D34').Value = '0.130'
$ws.Range('E34').Value = '  -0.37%  '
$ws.Range('B35').Value = 'Cosmos'
$ws.Range('C35').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('D35').Value = '13.30'
$ws.Range('E35').Value = '  -2.88%  '
$ws.Range('D36').Value = '682.73'
$ws.Range('E36').Value = '  +1.06%  '
$ws.Range('D37').Value = '65.03'
$ws.Range('E37').Value = '  -2.56%  '
$ws.Range('D38').Value = '0.434'
$ws.Range('E38').Value = '  +2.16%  '
$ws.Range('D39').Value = '0.0₃0814'
$ws.Range('E39').Value = '  -5.54%  '
$ws.Range('D40').Value = '0.149'
$ws.Range('E40').Value = '  -1.70%  '
$ws.Range('D41').Value = '3.38'
$ws.Range('E41').Value = '  -3.85%  '
$ws.Range('E43').Value = '  +0.39%  '
$ws.Range('D44').Value = '0.0479'
$ws.Range('E44').Value = '  -4.29%  '
$ws.Range('D45').Value = '3.19'
$ws.Range('E45').Value = '  -2.20%  '
$ws.Range('D46').Value = '0.148'
$ws.Range('E46').Value = '  -8.43%  '
$ws.Range('D47').Value = '2.67'
$ws.Range('E47').Value = '  -1.57%  '
$ws.Range('E48').Value = '  +5.25%  '
$ws.Range('D49').Value = '3.33'
$ws.Range('E49').Value = '  -5.51%  '
$ws.Range('D50').Value = '2.97'
$ws.Range('E50').Value = '  -2.75%  '
$ws.Range('D51').Value = '0.000269'
$ws.Range('E51').Value = '  -0.36%  '
